# 2020 update - part 1
# Import new data (new countries with their data quality score), appended
# after the existing table rows, and keep the table/range in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tableau1")

# New rows to append: ISO code + quality score.
$newData = @(
    @("BD", 1),
    @("BT", 1),
    @("KG", 1),
    @("KZ", 1),
    @("LA", 1),
    @("LK", 1),
    @("MM", 1),
    @("MN", 1),
    @("MV", 1),
    @("NP", 1),
    @("PH", 1),
    @("PK", 1),
    @("TJ", 1),
    @("TL", 1),
    @("TM", 1),
    @("UZ", 1),
    @("VN", 1),
    @("JP", 3)
)

foreach ($pair in $newData) {
    $newRow = $lo.ListRows.Add()
    $newRow.Range.Cells.Item(1, 1).Value = $pair[0]
    $newRow.Range.Cells.Item(1, 2).Value = $pair[1]
}

# Put the view where it ended up after the edits (bottom of the new data).
$lastCell = $lo.Range.Cells.Item($lo.Range.Rows.Count, 2)
$lastCell.Select()
$excel.ActiveWindow.ScrollRow = 206

$wb.Save()
